# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strike#) values for rows 2-15, column G
$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 5
    6  = 1
    7  = 5
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 0
    15 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
